# Apply the data update described by the commit "Fruta / hortaliza, semanal":
# three new weekly price records are inserted at the top of the existing block
# (pushing the rest of the block down by three rows), for
# Feria Lagunitas de Puerto Montt - Membrillo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows right before row 33; this shifts existing rows
# 33:83 down to 36:86 (and copies row formatting, e.g. the date style on
# column D, from the row above).
$ws.Rows("33:35").Insert()

# --- New row 33 ---------------------------------------------------------
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44665
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100104
$ws.Cells.Item(33, 8).Value = "Frutos de pepita"
$ws.Cells.Item(33, 9).Value = 100104003
$ws.Cells.Item(33, 10).Value = "Membrillo"
$ws.Cells.Item(33, 11).Value = "Champion"
$ws.Cells.Item(33, 12).Value = "Especial"
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(33, 14).Value = 18000
$ws.Cells.Item(33, 15).Value = 18000
$ws.Cells.Item(33, 16).Value = 18000
$ws.Cells.Item(33, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(33, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(33, 19).Value = 1000
$ws.Cells.Item(33, 20).Value = 18

# --- New row 34 ---------------------------------------------------------
$ws.Cells.Item(34, 1).Value = 4
$ws.Cells.Item(34, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value = "Los Lagos"
$ws.Cells.Item(34, 4).Value = 44665
$ws.Cells.Item(34, 5).Value = 10
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100104
$ws.Cells.Item(34, 8).Value = "Frutos de pepita"
$ws.Cells.Item(34, 9).Value = 100104003
$ws.Cells.Item(34, 10).Value = "Membrillo"
$ws.Cells.Item(34, 11).Value = "Champion"
$ws.Cells.Item(34, 12).Value = "Primera"
$ws.Cells.Item(34, 13).Value = 200
$ws.Cells.Item(34, 14).Value = 16000
$ws.Cells.Item(34, 15).Value = 16000
$ws.Cells.Item(34, 16).Value = 16000
$ws.Cells.Item(34, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(34, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(34, 19).Value = 889
$ws.Cells.Item(34, 20).Value = 18

# --- New row 35 ---------------------------------------------------------
$ws.Cells.Item(35, 1).Value = 4
$ws.Cells.Item(35, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(35, 3).Value = "Los Lagos"
$ws.Cells.Item(35, 4).Value = 44665
$ws.Cells.Item(35, 5).Value = 10
$ws.Cells.Item(35, 6).Value = "Fruta"
$ws.Cells.Item(35, 7).Value = 100104
$ws.Cells.Item(35, 8).Value = "Frutos de pepita"
$ws.Cells.Item(35, 9).Value = 100104003
$ws.Cells.Item(35, 10).Value = "Membrillo"
$ws.Cells.Item(35, 11).Value = "Champion"
$ws.Cells.Item(35, 12).Value = "Segunda"
$ws.Cells.Item(35, 13).Value = 200
$ws.Cells.Item(35, 14).Value = 13000
$ws.Cells.Item(35, 15).Value = 13000
$ws.Cells.Item(35, 16).Value = 13000
$ws.Cells.Item(35, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(35, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(35, 19).Value = 722
$ws.Cells.Item(35, 20).Value = 18
